$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.969.23'
$ws.Range("E2").Value = '  +2.57%  '

$ws.Range("D3").Value = '1.679.03'
$ws.Range("E3").Value = '  +1.96%  '

$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = '  -0.57%  '

$ws.Range("D5").Value = "'328.94"
$ws.Range("E5").Value = '  +7.44%  '

$ws.Range("D6").Value = "'0.9976"
$ws.Range("E6").Value = '  -0.33%  '

$ws.Range("E7").Value = '  +0.90%  '

$ws.Range("D8").Value = "'47.18"
$ws.Range("E8").Value = '  -0.29%  '

$ws.Range("D9").Value = "'0.3265"
$ws.Range("E9").Value = '  -0.13%  '

$ws.Range("D10").Value = "'1.151"
$ws.Range("E10").Value = '  +3.70%  '

$ws.Range("D11").Value = "'0.07107"
$ws.Range("E11").Value = '  +3.24%  '

$ws.Range("D12").Value = "'0.9970"
$ws.Range("E12").Value = '  -0.57%  '

$ws.Range("D13").Value = "'6.117"
$ws.Range("E13").Value = '  +3.38%  '

$ws.Range("D14").Value = "'19.81"
$ws.Range("E14").Value = '  +4.22%  '

$ws.Range("D15").Value = '1.675.90'
$ws.Range("E15").Value = '  +1.89%  '

$ws.Range("D16").Value = "'6.645"
$ws.Range("E16").Value = '  +1.88%  '

$ws.Range("D17").Value = "'0.00001056"
$ws.Range("E17").Value = '  +1.78%  '

$ws.Range("E18").Value = '  +1.78%  '

$ws.Range("E19").Value = '  -0.30%  '

$ws.Range("D20").Value = "'79.26"
$ws.Range("E20").Value = '  +4.00%  '

$ws.Range("D21").Value = "'16.01"
$ws.Range("E21").Value = '  +2.99%  '

$ws.Range("D22").Value = "'5.958"
$ws.Range("E22").Value = '  +1.41%  '

$ws.Range("D23").Value = "'12.75"
$ws.Range("E23").Value = '  +5.20%  '

$ws.Range("D24").Value = '24.957.68'
$ws.Range("E24").Value = '  +2.53%  '

$ws.Range("D25").Value = "'2.458"
$ws.Range("E25").Value = '  +1.15%  '

$ws.Range("D26").Value = "'2.431"
$ws.Range("E26").Value = '  +6.47%  '

$ws.Range("D27").Value = "'148.79"
$ws.Range("E27").Value = '  +2.32%  '

$ws.Range("D28").Value = "'18.86"
$ws.Range("E28").Value = '  +2.87%  '

$ws.Range("D29").Value = '1.862.10'
$ws.Range("E29").Value = '  +1.70%  '

$ws.Range("D30").Value = "'126.33"
$ws.Range("E30").Value = '  +1.96%  '

$ws.Range("D31").Value = "'1.194"
$ws.Range("E31").Value = '  +4.17%  '

$ws.Range("D32").Value = "'4.080"
$ws.Range("E32").Value = '  +0.97%  '

$ws.Range("D33").Value = "'5.799"
$ws.Range("E33").Value = '  +5.15%  '

$ws.Range("D34").Value = "'0.08472"
$ws.Range("E34").Value = '  +1.99%  '

$ws.Range("D35").Value = "'1.648"
$ws.Range("E35").Value = '  -1.47%  '

$ws.Range("D36").Value = "'12.34"
$ws.Range("E36").Value = '  +2.07%  '

$ws.Range("D37").Value = "'5.198"
$ws.Range("E37").Value = '  +0.69%  '

$ws.Range("D38").Value = "'0.02277"
$ws.Range("E38").Value = '  +3.58%  '

$ws.Range("D39").Value = "'0.06083"
$ws.Range("E39").Value = '  +1.42%  '

$ws.Range("E40").Value = '  +2.75%  '

$ws.Range("D41").Value = "'0.2095"
$ws.Range("E41").Value = '  +3.29%  '

$ws.Range("D42").Value = "'8.313"
$ws.Range("E42").Value = '  +1.76%  '

$ws.Range("D43").Value = "'0.9978"
$ws.Range("E43").Value = '  -0.20%  '

$ws.Range("D44").Value = "'0.5986"
$ws.Range("E44").Value = '  +3.49%  '

$ws.Range("D45").Value = "'13.63"
$ws.Range("E45").Value = '  +8.05%  '

$ws.Range("D46").Value = "'3.845"
$ws.Range("E46").Value = '  +3.57%  '

$ws.Range("D47").Value = "'0.5748"
$ws.Range("E47").Value = '  +4.26%  '

$ws.Range("D48").Value = "'126.04"
$ws.Range("E48").Value = '  +4.35%  '

$ws.Range("D49").Value = "'1.970"
$ws.Range("E49").Value = '  +2.57%  '

$ws.Range("D50").Value = "'0.07030"
$ws.Range("E50").Value = '  +2.31%  '

$ws.Range("D51").Value = "'1.193"
$ws.Range("E51").Value = '  +4.54%  '
